# RPA datasets push 2024-04-24
# Remove the "SK" (SK증권제11호스팩) and "유안타" (유안타제15호스팩) underwriter
# rows from the form3 sheet; remaining rows shift up accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 holds the SK entry, row 10 holds the 유안타 entry (original layout).
# Delete the lower row first so row 6's index is unaffected by the second delete.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(6).Delete()
